$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows before the old row 3 (20170103 entry), shifting
# all the existing trade-log rows down by 3. This recreates the row
# layout: row3 (blank spacer), row4 (new trade entry), row5 (blank
# spacer), row6 (the 20170103 entry that used to be row3), etc.
$ws.Rows("3:5").Insert()

# Populate the new trade-plan entry for 20170104 / 20170105 in row 4.
$ws.Range("A4").Value = 20170104
$ws.Range("B4").Value = "Wednesday"
$ws.Range("C4").Value = 20170105
$ws.Range("D4").Value = "Thursday"
$ws.Range("E4").Value = "Not sure it's a great move to sold most of DGAZ today at 4.31 but it's good and conservative move, as there're always plenty of opportunities out there in the market. As for tomorrow the report day, the consensus is ~-72 to -85, but keep in mind that even the weather is now warmer than normal, the inventory is decreasing than last year and five-year average, so even normal or warmer weather will continue to use up gas in storage, thus this should prevent a total collapse in prices. The current support might be seen around `$3.1, and the dropping of NG is slowing down, we can see what will happen around `$3.1"
$ws.Range("F4").Value = "hold and watch, don't buy or sell before clear signal"
$ws.Range("G4").Value = 20170104

# Row 4 wraps onto 7 lines of text (matches the 100.8pt height used by the
# other wrapped, multi-line trade entries in this sheet).
$ws.Rows(4).RowHeight = 100.8

# Move / update the active selection to the new entry's Lessons cell.
$ws.Range("F4").Select() | Out-Null
